# Update the cosinor analysis results for the re-run CircadiPy simulation
# (square wave, 0.1 amplitude) — commit "Make figures again to publication".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("E2").Value = 25.97000000000062
$ws.Range("H2").Value = [double]"6.00120553851436e-16"
$ws.Range("K2").Value = 54.37936759898865
$ws.Range("L2").Value = "[44.33140403004295, 64.42733116793435]"
$ws.Range("O2").Value = 1.855395060678656
$ws.Range("P2").Value = "[1.6667108172198102, 2.0440793041375027]"
$ws.Range("S2").Value = 65.38490285827498
$ws.Range("T2").Value = "[59.585753274818494, 71.18405244173147]"
$ws.Range("W2").Value = 18.30118118118162
$ws.Range("X2").Value = 17.52130130130172
$ws.Range("Y2").Value = 19.08106106106152

# ---- Row 3 ----
$ws.Range("E3").Value = 25.15000000000049
$ws.Range("G3").Value = [double]"1.425641604768657e-09"
$ws.Range("H3").Value = [double]"4.808146638037282e-09"
$ws.Range("K3").Value = 44.54113572278699
$ws.Range("L3").Value = "[26.96111616588948, 62.1211552796845]"
$ws.Range("M3").Value = [double]"1.612199053147023e-06"
$ws.Range("N3").Value = [double]"1.612199053147023e-06"
$ws.Range("O3").Value = 0.9371317425122712
$ws.Range("P3").Value = "[0.5471843060306556, 1.3270791789938867]"
$ws.Range("Q3").Value = [double]"4.929529937136579e-06"
$ws.Range("R3").Value = [double]"4.929529937136579e-06"
$ws.Range("S3").Value = 63.30291238680987
$ws.Range("T3").Value = "[54.03561676678154, 72.57020800683821]"
$ws.Range("W3").Value = 21.39889889889931
$ws.Range("X3").Value = 19.83803803803843
$ws.Range("Y3").Value = 22.9597597597602
